# Update the "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly generated output (see commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 68
$wsExhibit.Range("F3").Value = 525
$wsExhibit.Range("F4").Value = 171
$wsExhibit.Range("F5").Value = 222
$wsExhibit.Range("F6").Value = 369
$wsExhibit.Range("F7").Value = 230
$wsExhibit.Range("F8").Value = 2240
$wsExhibit.Range("F10").Value = 5514
$wsExhibit.Range("F11").Value = 128

# --- Sheet: 全部类型 ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 68
$wsAll.Range("F4").Value = 525
$wsAll.Range("F5").Value = 171
$wsAll.Range("F6").Value = 222
$wsAll.Range("F7").Value = 369
$wsAll.Range("F8").Value = 230
$wsAll.Range("F11").Value = 2240
$wsAll.Range("F13").Value = 5514
$wsAll.Range("F14").Value = 128
